# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) was mis-derived as the literal file-name string
# "5-1-2013-14" for every row; it should be the real game date, formatted
# as "2014-05-01".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateCol = 58   # column BF
$oldText = "5-1-2013-14"
$newText = "2014-05-01"

# Find the last used row in the Date column (header is row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, $dateCol).End(-4162).Row

$rng = $ws.Range($ws.Cells.Item(2, $dateCol), $ws.Cells.Item($lastRow, $dateCol))

# Writing the literal string straight into .Value would make Excel's
# automatic data-type detection parse "2014-05-01" as a real date serial
# number (since it matches a date pattern). Going through a text formula
# and converting the formula results to static values keeps the cells as
# plain text, exactly like the original "5-1-2013-14" strings were.
$rng.Formula = '="' + $newText + '"'
$rng.Copy()
$rng.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
